$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "0.584", "74.80").
# Force text format first so Excel keeps the exact authored string (incl.
# trailing zeros) instead of silently converting to a float, then reset the
# number format back to General/Normal so no stray style lingers on the cell.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D15", "D16", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D30", "D31", "D32", "D33", "D34", "D38", "D40", "D41", "D43", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.411.72"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "2.237.80"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "318.60"
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("D6").Value = "100.66"
$ws.Range("E6").Value = "  +3.68%  "
$ws.Range("D7").Value = "0.584"
$ws.Range("E7").Value = "  +2.89%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").Value = "37.46"
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").Value = "0.0835"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").Value = "7.73"
$ws.Range("E12").Value = "  +3.38%  "
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.579.34"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.867"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "14.32"
$ws.Range("E16").Value = "  +4.42%  "
$ws.Range("D17").Value = "2.256.49"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "43.384.77"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").Value = "14.24"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("E20").Value = "  +5.96%  "
$ws.Range("D21").Value = "6.64"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").Value = "65.63"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").Value = "3.16"
$ws.Range("D24").Value = "236.95"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  +4.95%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +2.87%  "
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("D30").Value = "6.39"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "36.62"
$ws.Range("E31").Value = "  +11.52%  "
$ws.Range("D32").Value = "20.32"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "0.0874"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").Value = "159.88"
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "1.89"
$ws.Range("E38").Value = "  +5.26%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("D41").Value = "3.73"
$ws.Range("E41").Value = "  +8.58%  "
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("D43").Value = "14.85"
$ws.Range("E43").Value = "  +27.15%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "1.816.11"
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("D46").Value = "0.205"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").Value = "83.97"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").Value = "5.31"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "8.81"
$ws.Range("E49").Value = "  +4.79%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "74.80"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "59.02"
$ws.Range("E51").Value = "  -0.09%  "

foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
